$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 23
$ws.Range("I2").Value = 64
$ws.Range("J2").Value = 250
$ws.Range("L2").Value = 76
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 52
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 34
$ws.Range("T2").Value = 43
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 424
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 446
$ws.Range("Z2").Value = 5
